$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(87, 8).Value = 94784.664
$ws.Cells.Item(87, 10).Value = 95141.60000000001
$ws.Cells.Item(87, 12).Value = 95141.60000000001
$ws.Cells.Item(87, 14).Value = -97637.60000000001

$ws.Cells.Item(90, 8).Value = 94784.664
$ws.Cells.Item(90, 10).Value = 95141.60000000001
$ws.Cells.Item(90, 12).Value = 285424.8
$ws.Cells.Item(90, 14).Value = -297904.8

$ws.Cells.Item(98, 8).Value = 2407.4285
$ws.Cells.Item(98, 9).Value = 1142
$ws.Cells.Item(98, 11).Value = 1142
$ws.Cells.Item(98, 13).Value = 356

$ws.Cells.Item(122, 8).Value = 2407.4285
$ws.Cells.Item(122, 9).Value = 1142
$ws.Cells.Item(122, 11).Value = 3426
$ws.Cells.Item(122, 13).Value = -976

$ws.Cells.Item(132, 8).Value = 1356.4546
$ws.Cells.Item(132, 9).Value = 1261.2941
$ws.Cells.Item(132, 11).Value = 3783.8823
$ws.Cells.Item(132, 13).Value = -1253.8823

$ws.Cells.Item(135, 8).Value = 1003.4
$ws.Cells.Item(135, 9).Value = 950.9474
$ws.Cells.Item(135, 11).Value = 8558.526600000001
$ws.Cells.Item(135, 13).Value = -6023.526600000001

$ws.Cells.Item(137, 8).Value = 2609.158
$ws.Cells.Item(137, 9).Value = 1972.5333
$ws.Cells.Item(137, 10).Value = 4996.5
$ws.Cells.Item(137, 11).Value = 5917.5999
$ws.Cells.Item(137, 12).Value = 14989.5
$ws.Cells.Item(137, 13).Value = -3367.5999
$ws.Cells.Item(137, 14).Value = -20089.5

$ws.Cells.Item(138, 8).Value = 4477.488
$ws.Cells.Item(138, 10).Value = 4868.0303
$ws.Cells.Item(138, 12).Value = 14604.0909
$ws.Cells.Item(138, 14).Value = -24884.0909

$ws.Cells.Item(141, 8).Value = 3667.647
$ws.Cells.Item(141, 9).Value = 3185.8572
$ws.Cells.Item(141, 11).Value = 9557.571599999999
$ws.Cells.Item(141, 13).Value = -4377.571599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3707.3333
$ws.Cells.Item(32, 9).Value = 3529.3333
$ws.Cells.Item(32, 11).Value = 3529.3333
$ws.Cells.Item(32, 13).Value = -3242.3333

$ws.Cells.Item(45, 8).Value = 3264.5386
$ws.Cells.Item(45, 9).Value = 2944.2
$ws.Cells.Item(45, 11).Value = 2944.2
$ws.Cells.Item(45, 13).Value = -2567.2

$ws.Cells.Item(74, 8).Value = 4348516
$ws.Cells.Item(74, 9).Value = 5128278.5
$ws.Cells.Item(74, 10).Value = 4125.7144
$ws.Cells.Item(74, 11).Value = 5128278.5
$ws.Cells.Item(74, 12).Value = 4125.7144
$ws.Cells.Item(74, 13).Value = -5127404.5
$ws.Cells.Item(74, 14).Value = -5873.7144

$ws.Cells.Item(77, 8).Value = 4348516
$ws.Cells.Item(77, 9).Value = 5128278.5
$ws.Cells.Item(77, 10).Value = 4125.7144
$ws.Cells.Item(77, 11).Value = 25641392.5
$ws.Cells.Item(77, 12).Value = 20628.572
$ws.Cells.Item(77, 13).Value = -25637024.5
$ws.Cells.Item(77, 14).Value = -29364.572

$ws.Cells.Item(80, 8).Value = 78110
$ws.Cells.Item(80, 10).Value = 78110
$ws.Cells.Item(80, 12).Value = 78110
$ws.Cells.Item(80, 14).Value = -80106

$ws.Cells.Item(83, 8).Value = 78110
$ws.Cells.Item(83, 10).Value = 78110
$ws.Cells.Item(83, 12).Value = 234330
$ws.Cells.Item(83, 14).Value = -244314

$ws.Cells.Item(122, 8).Value = 3213.4285
$ws.Cells.Item(122, 9).Value = 3213.4285
$ws.Cells.Item(122, 11).Value = 9640.2855
$ws.Cells.Item(122, 13).Value = -7190.2855

$ws.Cells.Item(132, 8).Value = 2152.0244
$ws.Cells.Item(132, 9).Value = 1038.5555
$ws.Cells.Item(132, 11).Value = 3115.6665
$ws.Cells.Item(132, 13).Value = -585.6664999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2724.9678
$ws.Cells.Item(134, 9).Value = 2763.9666
$ws.Cells.Item(134, 11).Value = 8291.899800000001
$ws.Cells.Item(134, 13).Value = -5756.899800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 157.90909
$ws.Cells.Item(7, 10).Value = 224.25
$ws.Cells.Item(7, 12).Value = 224.25
$ws.Cells.Item(7, 14).Value = -450.25

$ws.Cells.Item(31, 8).Value = 2461.1082
$ws.Cells.Item(31, 9).Value = 2354.3823
$ws.Cells.Item(31, 10).Value = 3670.6667
$ws.Cells.Item(31, 11).Value = 2354.3823
$ws.Cells.Item(31, 12).Value = 3670.6667
$ws.Cells.Item(31, 13).Value = -2059.3823
$ws.Cells.Item(31, 14).Value = -4260.6667

$ws.Cells.Item(34, 8).Value = 2461.1082
$ws.Cells.Item(34, 9).Value = 2354.3823
$ws.Cells.Item(34, 10).Value = 3670.6667
$ws.Cells.Item(34, 11).Value = 2354.3823
$ws.Cells.Item(34, 12).Value = 3670.6667
$ws.Cells.Item(34, 13).Value = -2152.3823
$ws.Cells.Item(34, 14).Value = -4074.6667

$ws.Cells.Item(58, 8).Value = 3013.25
$ws.Cells.Item(58, 9).Value = 3012
$ws.Cells.Item(58, 11).Value = 3012
$ws.Cells.Item(58, 13).Value = -2809

$ws.Cells.Item(86, 8).Value = 9002.200000000001
$ws.Cells.Item(86, 10).Value = 6669.5
$ws.Cells.Item(86, 12).Value = 6669.5
$ws.Cells.Item(86, 14).Value = -8915.5

$ws.Cells.Item(89, 8).Value = 9002.200000000001
$ws.Cells.Item(89, 10).Value = 6669.5
$ws.Cells.Item(89, 12).Value = 33347.5
$ws.Cells.Item(89, 14).Value = -44579.5

$ws.Cells.Item(99, 8).Value = 4327.25
$ws.Cells.Item(99, 9).Value = 4327.25
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 4327.25
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -2829.25
$ws.Cells.Item(99, 14).ClearContents()

$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 9).Value = 0
$ws.Cells.Item(105, 11).Value = 0
$ws.Cells.Item(105, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 3049.5557
$ws.Cells.Item(122, 9).Value = 2741
$ws.Cells.Item(122, 10).Value = 3666.6667
$ws.Cells.Item(122, 11).Value = 8223
$ws.Cells.Item(122, 12).Value = 11000.0001
$ws.Cells.Item(122, 13).Value = -5773
$ws.Cells.Item(122, 14).Value = -15900.0001

$ws.Cells.Item(126, 8).Value = 4327.25
$ws.Cells.Item(126, 9).Value = 4327.25
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 12981.75
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -10511.75
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 3361
$ws.Cells.Item(132, 9).Value = 2269
$ws.Cells.Item(132, 10).Value = 4999
$ws.Cells.Item(132, 11).Value = 6807
$ws.Cells.Item(132, 12).Value = 14997
$ws.Cells.Item(132, 13).Value = -4277
$ws.Cells.Item(132, 14).Value = -20057

$ws.Cells.Item(136, 8).Value = 3013.25
$ws.Cells.Item(136, 9).Value = 3012
$ws.Cells.Item(136, 11).Value = 9036
$ws.Cells.Item(136, 13).Value = -6486

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1907.3334
$ws.Cells.Item(97, 9).Value = 2205.8
$ws.Cells.Item(97, 10).Value = 415
$ws.Cells.Item(97, 11).Value = 2205.8
$ws.Cells.Item(97, 12).Value = 415
$ws.Cells.Item(97, 13).Value = -1709.8
$ws.Cells.Item(97, 14).Value = -1407

$ws.Cells.Item(122, 8).Value = 1942.8
$ws.Cells.Item(122, 9).Value = 1633
$ws.Cells.Item(122, 11).Value = 4899
$ws.Cells.Item(122, 13).Value = -2449

$ws.Cells.Item(132, 8).Value = 1109.2572
$ws.Cells.Item(132, 9).Value = 671.8387
$ws.Cells.Item(132, 11).Value = 2015.5161
$ws.Cells.Item(132, 13).Value = 514.4838999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 534.1667
$ws.Cells.Item(22, 9).Value = 415.66666
$ws.Cells.Item(22, 10).Value = 652.6667
$ws.Cells.Item(22, 11).Value = 415.66666
$ws.Cells.Item(22, 12).Value = 652.6667
$ws.Cells.Item(22, 13).Value = -120.66666
$ws.Cells.Item(22, 14).Value = -1242.6667

$ws.Cells.Item(27, 8).Value = 534.1667
$ws.Cells.Item(27, 9).Value = 415.66666
$ws.Cells.Item(27, 10).Value = 652.6667
$ws.Cells.Item(27, 11).Value = 415.66666
$ws.Cells.Item(27, 12).Value = 652.6667
$ws.Cells.Item(27, 13).Value = -308.66666
$ws.Cells.Item(27, 14).Value = -866.6667

$ws.Cells.Item(46, 8).Value = 171.42857

$ws.Cells.Item(68, 8).Value = 3000
$ws.Cells.Item(68, 9).Value = 3000
$ws.Cells.Item(68, 11).Value = 3000
$ws.Cells.Item(68, 13).Value = -2251

$ws.Cells.Item(71, 8).Value = 3000
$ws.Cells.Item(71, 9).Value = 3000
$ws.Cells.Item(71, 11).Value = 15000
$ws.Cells.Item(71, 13).Value = -11256

$ws.Cells.Item(132, 8).Value = 2458.8064
$ws.Cells.Item(132, 9).Value = 1535.1904
$ws.Cells.Item(132, 11).Value = 4605.5712
$ws.Cells.Item(132, 13).Value = -2075.5712

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1451.6857
$ws.Cells.Item(136, 9).Value = 1144.44
$ws.Cells.Item(136, 10).Value = 2219.8
$ws.Cells.Item(136, 11).Value = 3433.32
$ws.Cells.Item(136, 12).Value = 6659.400000000001
$ws.Cells.Item(136, 13).Value = -883.3200000000002
$ws.Cells.Item(136, 14).Value = -11759.4
